$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-entered dates for rows 6-11 (column B)
$ws.Range("B6").Value = 43866
$ws.Range("B7").Value = 43867
$ws.Range("B8").Value = 43867
$ws.Range("B9").Value = 43867
$ws.Range("B10").Value = 43867
$ws.Range("B11").Value = 43867

# Update the existing date in B12 from 06/02/2020 to 07/02/2020
$ws.Range("B12").Value = 43868

# Update the active selection to B11
$ws.Range("B11").Select()
